$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProcessAndPolicies")

# Insert a new row at position 6 (pushes the existing rows 6,7 down to 7,8)
$ws.Rows.Item(6).Insert()

# Bring over the formatting (column styles) from the row above (row 4), which
# has the same per-column style pattern as the new row needs, and from row 5
# for columns A:B (the "NullValue" placeholder style used by every row below
# the first data row).
$ws.Range("C4:N4").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("A5:B5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New test case data (fetch/Create verb row) for row 6
$ws.Range("A6").Value = "NullValue"
$ws.Range("B6").Value = "NullValue"
$ws.Range("C6").Value = "Create"
$ws.Range("D6").Value = "formTemplate"
$ws.Range("E6").Value = "EmployeeDetails"
$ws.Range("F6").Value = "employeeName"
$ws.Range("G6").Value = "any"
$ws.Range("H6").Value = "EmployeeDetails.employeeName,not equals,KGM100 and`nEmployeeDetails.employeeName,not equals,KGM120"
$ws.Range("K6").Value = "EmployeeDetails.employeeName,not equals,KGM211 and`nEmployeeDetails.employeeName,not equals,KGM221"

$ws.Rows.Item(6).RowHeight = 30

# Update the active selection to reflect where the edit left the cursor
$ws.Range("I5").Select()
